# Updates cryptos list values (Price column D and Volume(1h) column E)
# to reflect the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "76.485.95"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.085.75"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +4.13%  "
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "199.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "621.59"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.18%  "
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.216"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +5.47%  "
$ws.Range("E9").Value = "  -0.09%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.459"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("E11").Value = "  +0.12%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.24"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +6.72%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.652.07"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.09%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "29.66"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.79%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000202"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.85%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "76.417.89"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.073.42"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.01%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.57"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "9.16"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.57%  "
$ws.Range("E20").Value = "  +19.32%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "386.14"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.94%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.54"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.62%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +6.55%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.231.44"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.46%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "72.83"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.35"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.96%  "
$ws.Range("E28").Value = "  -0.07%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.0000112"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  +4.79%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "509.76"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("E34").Value = "  +6.26%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.132"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +18.41%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "20.94"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.99%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "163.21"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.94%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "195.46"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +8.20%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  -3.97%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -5.62%  "
$ws.Range("E43").Value = "  +0.07%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.28"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.65%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.796"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +19.06%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.28"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.07%  "
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("E48").Value = "  +7.12%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "41.17"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("E50").Value = "  +1.56%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.96"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.99%  "
